$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '21.664.66'
Set-TextValue $ws 'E2' '  -1.58%  '
Set-TextValue $ws 'D3' '1.534.02'
Set-TextValue $ws 'E3' '  -1.41%  '
Set-TextValue $ws 'D4' '1.000'
Set-TextValue $ws 'E4' '  +0.06%  '
Set-TextValue $ws 'E5' '  +0.10%  '
Set-TextValue $ws 'D6' '288.43'
Set-TextValue $ws 'E6' '  +0.50%  '
Set-TextValue $ws 'D7' '0.3942'
Set-TextValue $ws 'E7' '  +2.69%  '
Set-TextValue $ws 'D8' '0.3150'
Set-TextValue $ws 'E8' '  -2.97%  '
Set-TextValue $ws 'D9' '42.28'
Set-TextValue $ws 'E9' '  +1.85%  '
Set-TextValue $ws 'D10' '0.07154'
Set-TextValue $ws 'E10' '  -2.22%  '
Set-TextValue $ws 'E11' '  -6.82%  '
Set-TextValue $ws 'E12' '  +0.05%  '
Set-TextValue $ws 'D13' '5.625'
Set-TextValue $ws 'E13' '  -1.42%  '
Set-TextValue $ws 'D14' '18.51'
Set-TextValue $ws 'E14' '  -4.16%  '
Set-TextValue $ws 'D15' '6.592'
Set-TextValue $ws 'E15' '  -2.95%  '
Set-TextValue $ws 'D16' '1.536.12'
Set-TextValue $ws 'E16' '  -0.96%  '
Set-TextValue $ws 'D17' '0.00001089'
Set-TextValue $ws 'E17' '  -0.16%  '
Set-TextValue $ws 'D18' '0.06585'
Set-TextValue $ws 'E18' '  -0.62%  '
Set-TextValue $ws 'D19' '83.03'
Set-TextValue $ws 'E19' '  -2.51%  '
Set-TextValue $ws 'E20' '  +0.11%  '
Set-TextValue $ws 'D21' '6.096'
Set-TextValue $ws 'E21' '  -4.68%  '
Set-TextValue $ws 'D22' '15.38'
Set-TextValue $ws 'E22' '  -3.34%  '
Set-TextValue $ws 'D23' '10.85'
Set-TextValue $ws 'E23' '  -4.89%  '
Set-TextValue $ws 'D24' '2.383'
Set-TextValue $ws 'E24' '  +3.23%  '
Set-TextValue $ws 'D25' '21.671.20'
Set-TextValue $ws 'E25' '  -1.61%  '
Set-TextValue $ws 'D26' '2.338'
Set-TextValue $ws 'E26' '  -7.33%  '
Set-TextValue $ws 'D27' '147.93'
Set-TextValue $ws 'D28' '18.28'
Set-TextValue $ws 'E28' '  -3.08%  '
Set-TextValue $ws 'D29' '4.841'
Set-TextValue $ws 'E29' '  -0.42%  '
Set-TextValue $ws 'D30' '1.708.15'
Set-TextValue $ws 'E30' '  -1.13%  '
Set-TextValue $ws 'D31' '116.83'
Set-TextValue $ws 'E31' '  -3.19%  '
Set-TextValue $ws 'D32' '5.846'
Set-TextValue $ws 'E32' '  -0.45%  '
Set-TextValue $ws 'D33' '0.9399'
Set-TextValue $ws 'E33' '  -14.08%  '
Set-TextValue $ws 'D34' '0.08137'
Set-TextValue $ws 'E34' '  -0.20%  '
Set-TextValue $ws 'D35' '8.438'
Set-TextValue $ws 'E35' '  -8.66%  '
Set-TextValue $ws 'D36' '0.06045'
Set-TextValue $ws 'E36' '  -2.48%  '
Set-TextValue $ws 'D37' '5.086'
Set-TextValue $ws 'E37' '  -3.00%  '
Set-TextValue $ws 'D38' '0.02203'
Set-TextValue $ws 'E38' '  -4.07%  '
Set-TextValue $ws 'D39' '1.443'
Set-TextValue $ws 'E39' '  -12.99%  '
Set-TextValue $ws 'D40' '0.2012'
Set-TextValue $ws 'E40' '  -4.35%  '
Set-TextValue $ws 'D41' '1.174'
Set-TextValue $ws 'E41' '  -3.64%  '
Set-TextValue $ws 'E42' '  +0.08%  '
Set-TextValue $ws 'D43' '10.85'
Set-TextValue $ws 'E43' '  +0.00%  '
Set-TextValue $ws 'D44' '0.5733'
Set-TextValue $ws 'E44' '  -3.23%  '
Set-TextValue $ws 'D45' '3.719'
Set-TextValue $ws 'E45' '  -0.08%  '
Set-TextValue $ws 'D46' '12.87'
Set-TextValue $ws 'E46' '  -4.19%  '
Set-TextValue $ws 'D47' '0.5469'
Set-TextValue $ws 'E47' '  -4.56%  '
Set-TextValue $ws 'D48' '1.158'
Set-TextValue $ws 'E48' '  +0.35%  '
Set-TextValue $ws 'B49' 'Quant'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws 'D49' '115.64'
Set-TextValue $ws 'E49' '  -2.96%  '
Set-TextValue $ws 'B50' 'NEARProtocol'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D50' '1.857'
Set-TextValue $ws 'E50' '  -3.78%  '
Set-TextValue $ws 'D51' '0.06685'
Set-TextValue $ws 'E51' '  -2.76%  '
